$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# Values that look numeric are forced to stay text (matching the original
# inlineStr cell type) by stamping a Text number format before the write
# and clearing formats afterward so no stray style survives.

$ws.Range("D2").Value = '67.332.41'
$ws.Range("E2").Value = '  +0.14%  '

$ws.Range("D3").Value = '2.629.39'
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.53'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.84'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.93%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.557'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.80%  '

$ws.Range("D9").Value = '2.629.14'
$ws.Range("E9").Value = '  +0.35%  '

$ws.Range("E10").Value = '  +3.48%  '

$ws.Range("E11").Value = '  +0.50%  '

$ws.Range("E12").Value = '  -0.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.351'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.71'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.26%  '

$ws.Range("D15").Value = '3.113.78'
$ws.Range("E15").Value = '  +0.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000182'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.61%  '

$ws.Range("D17").Value = '67.279.99'
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").Value = '2.633.07'
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.20'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '363.27'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.54'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.19%  '

$ws.Range("E22").Value = '  -0.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.12'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.15'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '66.33'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -7.42%  '

$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000102'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.09%  '

$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '577.20'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -7.82%  '

$ws.Range("E31").Value = '  -5.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.83'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.84'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.64%  '

$ws.Range("E34").Value = '  -3.86%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.52'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.92'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.15'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.37'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.369'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.27'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.60'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.23%  '

$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("E46").Value = '  -0.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '155.72'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.22%  '

$ws.Range("D48").Value = '0.0₆0288'
$ws.Range("E48").Value = '  -3.57%  '

$ws.Range("E49").Value = '  -1.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.623'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.67'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.81%  '
